$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; everything (old A..W) shifts right to (B..X).
$ws.Columns("A").Insert()

# New "Match ID" column: header label + constant value (30) for every player row.
$ws.Range("A2").Value = "Match ID"
$ws.Range("A4:A19").Value = 30

# Match the bold/no-border look of the header column (mirrors the neighbouring
# "Player ID" header cell, but without its border).
$rng = $ws.Range("A2:A19")
$rng.Font.Bold = $true

# Row 20 is the hidden "totals" row; fill it in too, toggling Hidden off/on
# around the write (and an AutoFit) so the engine doesn't stamp a stray
# custom row-height on the hidden row.
$r20 = $ws.Rows(20)
$r20.Hidden = $false
$ws.Range("A20").Value = 30
$r20.EntireRow.AutoFit()
$r20.Hidden = $true

# Match the selection left behind by the edit.
$ws.Range("A2:A19").Select() | Out-Null
